$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.040.36'
$ws.Range('E2').Value = '  +3.02%  '
$ws.Range('D3').Value = '2.456.23'
$ws.Range('E3').Value = '  +2.14%  '
$ws.Range('D5').Value = '576.86'
$ws.Range('E5').Value = '  +1.55%  '
$ws.Range('D6').Value = '146.48'
$ws.Range('E6').Value = '  +3.10%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').Value = '0.540'
$ws.Range('E8').Value = '  +0.77%  '
$ws.Range('D9').Value = '2.456.74'
$ws.Range('E9').Value = '  +1.77%  '
$ws.Range('E10').Value = '  +2.58%  '
$ws.Range('E11').Value = '  +2.49%  '
$ws.Range('D12').Value = '5.28'
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').Value = '0.355'
$ws.Range('E13').Value = '  +2.46%  '
$ws.Range('D14').Value = '28.42'
$ws.Range('E14').Value = '  +7.44%  '
$ws.Range('E15').Value = '  +5.04%  '
$ws.Range('D16').Value = '2.898.62'
$ws.Range('D17').Value = '63.042.72'
$ws.Range('E17').Value = '  +3.63%  '
$ws.Range('D18').Value = '2.458.13'
$ws.Range('E18').Value = '  +1.71%  '
$ws.Range('D19').Value = '7.94'
$ws.Range('E19').Value = '  -2.50%  '
$ws.Range('D20').Value = '11.06'
$ws.Range('E20').Value = '  +3.44%  '
$ws.Range('D21').Value = '330.23'
$ws.Range('E21').Value = '  +1.94%  '
$ws.Range('E22').Value = '  +1.09%  '
$ws.Range('D23').Value = '2.13'
$ws.Range('E23').Value = '  +10.13%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '66.36'
$ws.Range('E25').Value = '  +1.78%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').Value = '1.25'
$ws.Range('E26').Value = '  +25.26%  '
$ws.Range('D27').Value = '656.28'
$ws.Range('E27').Value = '  +10.51%  '
$ws.Range('D28').Value = '8.62'
$ws.Range('E28').Value = '  +4.32%  '
$ws.Range('D29').Value = '0.0000101'
$ws.Range('E29').Value = '  +6.77%  '
$ws.Range('D31').Value = '8.20'
$ws.Range('E31').Value = '  +2.21%  '
$ws.Range('D32').Value = '1.44'
$ws.Range('E32').Value = '  +4.96%  '
$ws.Range('E33').Value = '  +3.89%  '
$ws.Range('E34').Value = '  +4.77%  '
$ws.Range('D35').Value = '0.0₆0396'
$ws.Range('E35').Value = '  +40.57%  '
$ws.Range('E36').Value = '  +1.68%  '
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.26%  '
$ws.Range('D38').Value = '4.80'
$ws.Range('E38').Value = '  +3.82%  '
$ws.Range('D39').Value = '5.55'
$ws.Range('E39').Value = '  +5.34%  '
$ws.Range('D40').Value = '0.374'
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('D41').Value = '18.85'
$ws.Range('E41').Value = '  +2.67%  '
$ws.Range('D42').Value = '151.77'
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').Value = '2.75'
$ws.Range('E43').Value = '  +8.66%  '
$ws.Range('D44').Value = '1.77'
$ws.Range('E44').Value = '  +4.43%  '
$ws.Range('D45').Value = '42.65'
$ws.Range('E45').Value = '  +2.34%  '
$ws.Range('E47').Value = '  +27.27%  '
$ws.Range('D48').Value = '146.79'
$ws.Range('E48').Value = '  +3.70%  '
$ws.Range('D49').Value = '3.63'
$ws.Range('E49').Value = '  +2.86%  '
$ws.Range('D50').Value = '20.66'
$ws.Range('E50').Value = '  +4.05%  '
$ws.Range('D51').Value = '0.607'
$ws.Range('E51').Value = '  +2.57%  '
